# Adding a new indicator
# Appends a new "Investment" row (Expected investment) to the indicators
# table on the active worksheet, mirroring the manual edit captured in the
# commit: new row 66 with Category/Indicator/Rationale/Formula/Source data,
# plus the accompanying view changes (zoom level + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

$ws.Range("A$row").Value = "Investment"
$ws.Range("B$row").Value = "Expected investment"
$ws.Range("C$row").Value = "Expected Investment reflects the anticipated financial commitment investors are willing to allocate toward industry projects."
$ws.Range("D$row").Value = "Forecasted financial return of project i  x  The probability of project being funded, based on market analysis, investor surveys, and financial feasibility"
$ws.Range("E$row").Value = "Market research reports, investor surveys and commitments, historical data on investment trends in similar sectors"

# The Formula cell picked up an (empty/default) explicit alignment record
# when it was touched in the UI, producing a new cellXf entry.
$ws.Range("D$row").WrapText = $false

# View state: the workbook was re-saved at a different zoom level with the
# new row's Formula cell selected.
$excel.ActiveWindow.Zoom = 107
$ws.Range("D$row").Select()
